$wb = $excel.ActiveWorkbook

# --- Update "Status" text across all sheets (shared string text change) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

function Set-HandbackSheet($SheetName, $XlfName, $HandbackDateTime) {
    $ws = $wb.Worksheets.Item($SheetName)
    $mdName1 = "54091fa3-53b8-4d5f-893d-5888e79d5274.md"
    $mdName2 = "fffff303ce9d-5883-4e8e-bcc3-6ec604bb405f.md"

    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # New "Latest Target File" (E) / "Latest Handback File" (F) columns for rows 2 & 3
    $ws.Range("E2").Value = $mdName1
    $ws.Range("F2").Value = $XlfName
    $ws.Range("E3").Value = $mdName1
    $ws.Range("F3").Value = $XlfName

    $ws.Range("E2").Font.Underline = 2
    $ws.Range("E2").Font.Color = 15570276
    $ws.Range("F2").Font.Underline = 2
    $ws.Range("F2").Font.Color = 15570276
    $ws.Range("E3").Font.Underline = 2
    $ws.Range("E3").Font.Color = 15570276
    $ws.Range("F3").Font.Underline = 2
    $ws.Range("F3").Font.Color = 15570276

    # "Latest Handback DateTime" (G) for rows 2 & 3
    $ws.Range("G2").Value = $HandbackDateTime
    $ws.Range("G3").Value = $HandbackDateTime
}

Set-HandbackSheet "zh-cn" "54091fa3-53b8-4d5f-893d-5888e79d5274.b0e32b1ef89c78a0e721748aec972d6eb0b7cccd.zh-cn.xlf" "2016-01-26 07:14:27"
Set-HandbackSheet "de-de" "54091fa3-53b8-4d5f-893d-5888e79d5274.b0e32b1ef89c78a0e721748aec972d6eb0b7cccd.de-de.xlf" "2016-01-26 07:14:53"
